$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab (Sheet1 -> Formation Energy)
$ws.Name = "Formation Energy"

# New header cells for the additional structure columns (C1:G1)
$ws.Range("C1").Value = "single doping in the subsurface"
$ws.Range("D1").Value = "overlayer"
$ws.Range("E1").Value = "island"
$ws.Range("F1").Value = "subsurface"
$ws.Range("G1").Value = "parallelogram"

# Fill in the new data columns (C:G) for the existing dopant rows (2-8)
# Row 2: Ni
$ws.Range("C2").Value = 0.1092539975000042
$ws.Range("D2").Value = 0.36753135638888895
$ws.Range("E2").Value = 0.22382353250000531
$ws.Range("F2").Value = 0.31063214194444516
$ws.Range("G2").Value = 0.2323341250000075
# Row 3: Co
$ws.Range("C3").Value = 0.56391413500003473
$ws.Range("D3").Value = 0.72377658944444911
$ws.Range("E3").Value = 0.64230312499999975
$ws.Range("F3").Value = 0.65973504166667263
$ws.Range("G3").Value = 0.60245007000000506
# Row 4: V
$ws.Range("C4").Value = -0.20281884499997904
$ws.Range("D4").Value = -0.2061410638888839
$ws.Range("E4").Value = -0.066620834999992162
$ws.Range("F4").Value = -0.27891608944444474
$ws.Range("G4").Value = -0.040901902499993703
# Row 5: Cr
$ws.Range("C5").Value = -0.13964059499999504
$ws.Range("D5").Value = 0.3022197727777779
$ws.Range("E5").Value = 0.1503021475000077
$ws.Range("F5").Value = 0.39943862388889273
$ws.Range("G5").Value = 0.20930033750000465
# Row 6: Mn
$ws.Range("C6").Value = -0.71314150637931206
$ws.Range("D6").Value = -0.0092570241570888202
$ws.Range("E6").Value = -0.14955614637930736
$ws.Range("F6").Value = 0.054509480287357408
$ws.Range("G6").Value = -0.099625411379303586
# Row 7: Fe
$ws.Range("C7").Value = 0.37280091500000712
$ws.Range("D7").Value = 0.47790880055555796
$ws.Range("E7").Value = 0.53287292750000237
$ws.Range("F7").Value = 0.42017660611111096
$ws.Range("G7").Value = 0.49975075000001201
# Row 8: Pt
$ws.Range("C8").Value = 0.22755464750003362
$ws.Range("D8").Value = 0.2846047808333374
$ws.Range("E8").Value = 0.082740877500006249
$ws.Range("F8").Value = 0.52999399083333487
$ws.Range("G8").Value = 0.12850357000001189

# New dopant rows (9-10): Ti and Hf, matching style of existing rows
# Row 9: Ti
$ws.Range("A9").Value = "Ti"
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = -1.0553793694444416
$ws.Range("E9").Value = -1.1420011749999892
$ws.Range("F9").Value = -1.0078864494444422
$ws.Range("G9").Value = -1.0807497274999909
$ws.Range("B9:C9").NumberFormat = "0.000"
# Row 10: Hf
$ws.Range("A10").Value = "Hf"
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = -0.76025394722222039
$ws.Range("E10").Value = -1.5250302849999908
$ws.Range("F10").Value = -0.65028787499999929
$ws.Range("G10").Value = -1.4776542124999903
$ws.Range("B10:C10").NumberFormat = "0.000"

# Restore the reported selection from the saved workbook
$ws.Range("E19").Select()
